$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Pad home_address nicknames (rows 2-10: home_address_1..9 -> home_address_01..09)
for ($i = 1; $i -le 9; $i++) {
    $row = $i + 1
    $padded = "{0:D2}" -f $i
    $ws.Cells.Item($row, 2).Value = "home_address_$padded"
}

# 2. Pad seasonal_address nicknames (rows 13-21: seasonal_address_1..9 -> seasonal_address_01..09)
for ($i = 1; $i -le 9; $i++) {
    $row = $i + 12
    $padded = "{0:D2}" -f $i
    $ws.Cells.Item($row, 2).Value = "seasonal_address_$padded"
}

# 3. Insert a new row at row 23 (shifts old rows 23-25 down to 24-26)
$ws.Rows.Item(23).Insert()

# 4. Fill the newly inserted row 23 with the childhood_address_01 entry
# Force the cells to be stored as text (matching the other rows' inline string
# representation) so purely-numeric values like "632533534" don't turn into numbers.
$newRowValues = @(
    "632533534",
    "childhood_address_01",
    "284547539",
    "802585033",
    "746533238",
    "128827522",
    "439447560",
    "286781627",
    "733929451",
    "264797252",
    "890792569",
    "451394598",
    "984908796",
    "847327251",
    "469914719",
    "952124199",
    "204186397",
    "Childhood"
)
$newRowRange = $ws.Range($ws.Cells.Item(23, 1), $ws.Cells.Item(23, $newRowValues.Length))
$newRowRange.NumberFormat = "@"
for ($c = 1; $c -le $newRowValues.Length; $c++) {
    $ws.Cells.Item(23, $c).Value = $newRowValues[$c - 1]
}
# Drop the temporary text number-format now that the values are locked in as
# text, so the new row keeps the same (default) styling as the other rows.
$newRowRange.ClearFormats()

# 5. Pad nicknames for the rows that shifted down
$ws.Cells.Item(24, 2).Value = "current_work_address_01"
$ws.Cells.Item(25, 2).Value = "previous_work_address_01"
$ws.Cells.Item(26, 2).Value = "school_address_01"
